$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 40 (G=5505)
$ws_ALC.Cells.Item(40, 8).Value = 130562.125  # H40: 129562.375 -> 130562.125
$ws_ALC.Cells.Item(40, 10).Value = 148856.72  # J40: 147714.14 -> 148856.72
$ws_ALC.Cells.Item(40, 12).Value = 148856.72  # L40: 147714.14 -> 148856.72
$ws_ALC.Cells.Item(40, 14).Value = -149206.72  # N40: -148064.14 -> -149206.72

# ALC row 58 (G=4606)
$ws_ALC.Cells.Item(58, 8).Value = 38461944  # H58: 50000228 -> 38461944
$ws_ALC.Cells.Item(58, 9).Value = 161.66667  # I58: 171.25 -> 161.66667
$ws_ALC.Cells.Item(58, 10).Value = 125000950  # J58: 250000450 -> 125000950
$ws_ALC.Cells.Item(58, 11).Value = 485.00001  # K58: 513.75 -> 485.00001
$ws_ALC.Cells.Item(58, 12).Value = 375002850  # L58: 750001350 -> 375002850
$ws_ALC.Cells.Item(58, 13).Value = -335.00001  # M58: -363.75 -> -335.00001
$ws_ALC.Cells.Item(58, 14).Value = -375003150  # N58: -750001650 -> -375003150

# ALC row 98 (G=36237)
$ws_ALC.Cells.Item(98, 8).Value = 2122.6572  # H98: 2246.818 -> 2122.6572
$ws_ALC.Cells.Item(98, 9).Value = 2159.3438  # I98: 2298.3667 -> 2159.3438
$ws_ALC.Cells.Item(98, 11).Value = 2159.3438  # K98: 2298.3667 -> 2159.3438
$ws_ALC.Cells.Item(98, 13).Value = -661.3438000000001  # M98: -800.3667 -> -661.3438000000001

# ALC row 113 (G=27775)
$ws_ALC.Cells.Item(113, 8).Value = 71441450  # H113: 75013380 -> 71441450
$ws_ALC.Cells.Item(113, 9).Value = 4027.1667  # I113: 4272.6 -> 4027.1667
$ws_ALC.Cells.Item(113, 11).Value = 4027.1667  # K113: 4272.6 -> 4027.1667
$ws_ALC.Cells.Item(113, 13).Value = -773.1667000000002  # M113: -1018.6 -> -773.1667000000002

# ALC row 116 (G=27778)
$ws_ALC.Cells.Item(116, 8).Value = 13168793  # H116: 15635583 -> 13168793
$ws_ALC.Cells.Item(116, 9).Value = 20841674  # I116: 25008234 -> 20841674
$ws_ALC.Cells.Item(116, 10).Value = 15284.143  # J116: 14498.167 -> 15284.143
$ws_ALC.Cells.Item(116, 11).Value = 20841674  # K116: 25008234 -> 20841674
$ws_ALC.Cells.Item(116, 12).Value = 15284.143  # L116: 14498.167 -> 15284.143
$ws_ALC.Cells.Item(116, 13).Value = -20838232  # M116: -25004792 -> -20838232
$ws_ALC.Cells.Item(116, 14).Value = -22168.143  # N116: -21382.167 -> -22168.143

# ALC row 122 (G=36237)
$ws_ALC.Cells.Item(122, 8).Value = 2122.6572  # H122: 2246.818 -> 2122.6572
$ws_ALC.Cells.Item(122, 9).Value = 2159.3438  # I122: 2298.3667 -> 2159.3438
$ws_ALC.Cells.Item(122, 11).Value = 6478.0314  # K122: 6895.1001 -> 6478.0314
$ws_ALC.Cells.Item(122, 13).Value = -4028.0314  # M122: -4445.1001 -> -4028.0314

# ALC row 132 (G=44049)
$ws_ALC.Cells.Item(132, 8).Value = 2073.9443  # H132: 1248.9354 -> 2073.9443
$ws_ALC.Cells.Item(132, 9).Value = 1895.6875  # I132: 1093.6897 -> 1895.6875
$ws_ALC.Cells.Item(132, 11).Value = 5687.0625  # K132: 3281.0691 -> 5687.0625
$ws_ALC.Cells.Item(132, 13).Value = -3157.0625  # M132: -751.0690999999997 -> -3157.0625

# ALC row 141 (G=44161)
$ws_ALC.Cells.Item(141, 8).Value = 2731.9688  # H141: 2791.0967 -> 2731.9688
$ws_ALC.Cells.Item(141, 9).Value = 2152.6897  # I141: 2197.4644 -> 2152.6897
$ws_ALC.Cells.Item(141, 11).Value = 6458.0691  # K141: 6592.3932 -> 6458.0691
$ws_ALC.Cells.Item(141, 13).Value = -1278.0691  # M141: -1412.3932 -> -1278.0691

# ARM row 4 (G=5071)
$ws_ARM.Cells.Item(4, 8).Value = 200  # H4: 0 -> 200
$ws_ARM.Cells.Item(4, 10).Value = 200  # J4: 0 -> 200
$ws_ARM.Cells.Item(4, 12).Value = 200  # L4: 0 -> 200
$ws_ARM.Cells.Item(4, 14).Value = -432  # N4: None -> -432

# ARM row 5 (G=5091)
$ws_ARM.Cells.Item(5, 8).Value = 498.5  # H5: 0 -> 498.5
$ws_ARM.Cells.Item(5, 9).Value = 498  # I5: 0 -> 498
$ws_ARM.Cells.Item(5, 10).Value = 499  # J5: 0 -> 499
$ws_ARM.Cells.Item(5, 11).Value = 498  # K5: 0 -> 498
$ws_ARM.Cells.Item(5, 12).Value = 499  # L5: 0 -> 499
$ws_ARM.Cells.Item(5, 13).Value = -386  # M5: None -> -386
$ws_ARM.Cells.Item(5, 14).Value = -723  # N5: None -> -723

# ARM row 45 (G=27714)
$ws_ARM.Cells.Item(45, 8).Value = 3093.818  # H45: 3123.3 -> 3093.818
$ws_ARM.Cells.Item(45, 9).Value = 3198.3333  # I45: 3278.2 -> 3198.3333
$ws_ARM.Cells.Item(45, 11).Value = 3198.3333  # K45: 3278.2 -> 3198.3333
$ws_ARM.Cells.Item(45, 13).Value = -2821.3333  # M45: -2901.2 -> -2821.3333

# ARM row 61 (G=43999)
$ws_ARM.Cells.Item(61, 8).Value = 23811690  # H61: 25643262 -> 23811690
$ws_ARM.Cells.Item(61, 9).Value = 942.65625  # I61: 979.73334 -> 942.65625
$ws_ARM.Cells.Item(61, 10).Value = 100006080  # J61: 111117540 -> 100006080
$ws_ARM.Cells.Item(61, 11).Value = 942.65625  # K61: 979.73334 -> 942.65625
$ws_ARM.Cells.Item(61, 12).Value = 100006080  # L61: 111117540 -> 100006080
$ws_ARM.Cells.Item(61, 13).Value = -730.65625  # M61: -767.73334 -> -730.65625
$ws_ARM.Cells.Item(61, 14).Value = -100006504  # N61: -111117964 -> -100006504

# ARM row 74 (G=44000)
$ws_ARM.Cells.Item(74, 8).Value = 41277  # H74: 42888.12 -> 41277
$ws_ARM.Cells.Item(74, 9).Value = 57263.5  # I74: 60573.176 -> 57263.5
$ws_ARM.Cells.Item(74, 11).Value = 57263.5  # K74: 60573.176 -> 57263.5
$ws_ARM.Cells.Item(74, 13).Value = -56389.5  # M74: -59699.176 -> -56389.5

# ARM row 77 (G=44000)
$ws_ARM.Cells.Item(77, 8).Value = 41277  # H77: 42888.12 -> 41277
$ws_ARM.Cells.Item(77, 9).Value = 57263.5  # I77: 60573.176 -> 57263.5
$ws_ARM.Cells.Item(77, 11).Value = 286317.5  # K77: 302865.88 -> 286317.5
$ws_ARM.Cells.Item(77, 13).Value = -281949.5  # M77: -298497.88 -> -281949.5

# ARM row 132 (G=43997)
$ws_ARM.Cells.Item(132, 8).Value = 4934.451  # H132: 4652.7456 -> 4934.451
$ws_ARM.Cells.Item(132, 9).Value = 4321.8125  # I132: 3959.5 -> 4321.8125
$ws_ARM.Cells.Item(132, 11).Value = 12965.4375  # K132: 11878.5 -> 12965.4375
$ws_ARM.Cells.Item(132, 13).Value = -10435.4375  # M132: -9348.5 -> -10435.4375

# ARM row 136 (G=43999)
$ws_ARM.Cells.Item(136, 8).Value = 23811690  # H136: 25643262 -> 23811690
$ws_ARM.Cells.Item(136, 9).Value = 942.65625  # I136: 979.73334 -> 942.65625
$ws_ARM.Cells.Item(136, 10).Value = 100006080  # J136: 111117540 -> 100006080
$ws_ARM.Cells.Item(136, 11).Value = 2827.96875  # K136: 2939.20002 -> 2827.96875
$ws_ARM.Cells.Item(136, 12).Value = 300018240  # L136: 333352620 -> 300018240
$ws_ARM.Cells.Item(136, 13).Value = -277.96875  # M136: -389.2000200000002 -> -277.96875
$ws_ARM.Cells.Item(136, 14).Value = -300023340  # N136: -333357720 -> -300023340

# BSM row 4 (G=5091)
$ws_BSM.Cells.Item(4, 8).Value = 498.5  # H4: 0 -> 498.5
$ws_BSM.Cells.Item(4, 9).Value = 498  # I4: 0 -> 498
$ws_BSM.Cells.Item(4, 10).Value = 499  # J4: 0 -> 499
$ws_BSM.Cells.Item(4, 11).Value = 498  # K4: 0 -> 498
$ws_BSM.Cells.Item(4, 12).Value = 499  # L4: 0 -> 499
$ws_BSM.Cells.Item(4, 13).Value = -383  # M4: None -> -383
$ws_BSM.Cells.Item(4, 14).Value = -729  # N4: None -> -729

# BSM row 124 (G=34245)
$ws_BSM.Cells.Item(124, 8).Value = 0  # H124: 50567 -> 0
$ws_BSM.Cells.Item(124, 10).Value = 0  # J124: 50567 -> 0
$ws_BSM.Cells.Item(124, 12).Value = 0  # L124: 50567 -> 0
$ws_BSM.Cells.Item(124, 14).ClearContents()  # N124: -60387 -> (removed)

# BSM row 125 (G=34235)
$ws_BSM.Cells.Item(125, 8).Value = 0  # H125: 50052 -> 0
$ws_BSM.Cells.Item(125, 10).Value = 0  # J125: 50052 -> 0
$ws_BSM.Cells.Item(125, 12).Value = 0  # L125: 50052 -> 0
$ws_BSM.Cells.Item(125, 14).ClearContents()  # N125: -59892 -> (removed)

# BSM row 126 (G=34398)
$ws_BSM.Cells.Item(126, 8).Value = 0  # H126: 50567 -> 0
$ws_BSM.Cells.Item(126, 10).Value = 0  # J126: 50567 -> 0
$ws_BSM.Cells.Item(126, 12).Value = 0  # L126: 50567 -> 0
$ws_BSM.Cells.Item(126, 14).ClearContents()  # N126: -60447 -> (removed)

# BSM row 129 (G=35382)
$ws_BSM.Cells.Item(129, 8).Value = 49998  # H129: 51888.285 -> 49998
$ws_BSM.Cells.Item(129, 9).Value = 49998  # I129: 49999 -> 49998
$ws_BSM.Cells.Item(129, 10).Value = 0  # J129: 52203.168 -> 0
$ws_BSM.Cells.Item(129, 11).Value = 49998  # K129: 49999 -> 49998
$ws_BSM.Cells.Item(129, 12).Value = 0  # L129: 52203.168 -> 0
$ws_BSM.Cells.Item(129, 13).Value = -44998  # M129: -44999 -> -44998
$ws_BSM.Cells.Item(129, 14).ClearContents()  # N129: -62203.168 -> (removed)

# BSM row 130 (G=34682)
$ws_BSM.Cells.Item(130, 8).Value = 74240  # H130: 71193.10000000001 -> 74240
$ws_BSM.Cells.Item(130, 10).Value = 74240  # J130: 71193.10000000001 -> 74240
$ws_BSM.Cells.Item(130, 12).Value = 74240  # L130: 71193.10000000001 -> 74240
$ws_BSM.Cells.Item(130, 14).Value = -84280  # N130: -81233.10000000001 -> -84280

# BSM row 134 (G=43998)
$ws_BSM.Cells.Item(134, 8).Value = 7357144.5  # H134: 6253684.5 -> 7357144.5
$ws_BSM.Cells.Item(134, 9).Value = 11906772  # I134: 9260988 -> 11906772
$ws_BSM.Cells.Item(134, 11).Value = 35720316  # K134: 27782964 -> 35720316
$ws_BSM.Cells.Item(134, 13).Value = -35717781  # M134: -27780429 -> -35717781

# CRP row 62 (G=12580)
$ws_CRP.Cells.Item(62, 8).Value = 4858  # H62: 4879.905 -> 4858
$ws_CRP.Cells.Item(62, 9).Value = 4669.1763  # I62: 4686.125 -> 4669.1763
$ws_CRP.Cells.Item(62, 11).Value = 4669.1763  # K62: 4686.125 -> 4669.1763
$ws_CRP.Cells.Item(62, 13).Value = -4045.1763  # M62: -4062.125 -> -4045.1763

# CRP row 65 (G=12580)
$ws_CRP.Cells.Item(65, 8).Value = 4858  # H65: 4879.905 -> 4858
$ws_CRP.Cells.Item(65, 9).Value = 4669.1763  # I65: 4686.125 -> 4669.1763
$ws_CRP.Cells.Item(65, 11).Value = 23345.8815  # K65: 23430.625 -> 23345.8815
$ws_CRP.Cells.Item(65, 13).Value = -20225.8815  # M65: -20310.625 -> -20225.8815

# CRP row 132 (G=44019)
$ws_CRP.Cells.Item(132, 8).Value = 3268.2156  # H132: 3323.58 -> 3268.2156
$ws_CRP.Cells.Item(132, 9).Value = 2430.8918  # I132: 2484.5278 -> 2430.8918
$ws_CRP.Cells.Item(132, 11).Value = 7292.6754  # K132: 7453.5834 -> 7292.6754
$ws_CRP.Cells.Item(132, 13).Value = -4762.6754  # M132: -4923.5834 -> -4762.6754

# CUL row 38 (G=4860)
$ws_CUL.Cells.Item(38, 8).Value = 41666692  # H38: 35714316 -> 41666692
$ws_CUL.Cells.Item(38, 9).Value = 20.9  # I38: 28.09091 -> 20.9
$ws_CUL.Cells.Item(38, 10).Value = 250000050  # J38: 166666700 -> 250000050
$ws_CUL.Cells.Item(38, 11).Value = 62.7  # K38: 84.27273 -> 62.7
$ws_CUL.Cells.Item(38, 12).Value = 750000150  # L38: 500000100 -> 750000150
$ws_CUL.Cells.Item(38, 13).Value = 284.3  # M38: 262.72727 -> 284.3
$ws_CUL.Cells.Item(38, 14).Value = -750000844  # N38: -500000794 -> -750000844

# CUL row 103 (G=19839)
$ws_CUL.Cells.Item(103, 8).Value = 554  # H103: 37037532 -> 554
$ws_CUL.Cells.Item(103, 10).Value = 766.8  # J103: 55556196 -> 766.8
$ws_CUL.Cells.Item(103, 12).Value = 2300.4  # L103: 166668588 -> 2300.4
$ws_CUL.Cells.Item(103, 14).Value = -4058.4  # N103: -166670346 -> -4058.4

# GSM row 63 (G=11048)
$ws_GSM.Cells.Item(63, 8).Value = 0  # H63: 54999 -> 0
$ws_GSM.Cells.Item(63, 10).Value = 0  # J63: 54999 -> 0
$ws_GSM.Cells.Item(63, 12).Value = 0  # L63: 54999 -> 0
$ws_GSM.Cells.Item(63, 14).ClearContents()  # N63: -56371 -> (removed)

# GSM row 66 (G=11048)
$ws_GSM.Cells.Item(66, 8).Value = 0  # H66: 54999 -> 0
$ws_GSM.Cells.Item(66, 10).Value = 0  # J66: 54999 -> 0
$ws_GSM.Cells.Item(66, 12).Value = 0  # L66: 164997 -> 0
$ws_GSM.Cells.Item(66, 14).ClearContents()  # N66: -171861 -> (removed)

# GSM row 107 (G=27802)
$ws_GSM.Cells.Item(107, 8).Value = 471243.34  # H107: 693.375 -> 471243.34
$ws_GSM.Cells.Item(107, 9).Value = 1143056.2  # I107: 278.8 -> 1143056.2
$ws_GSM.Cells.Item(107, 10).Value = 974.3  # J107: 881.8182 -> 974.3
$ws_GSM.Cells.Item(107, 11).Value = 1143056.2  # K107: 278.8 -> 1143056.2
$ws_GSM.Cells.Item(107, 12).Value = 974.3  # L107: 881.8182 -> 974.3
$ws_GSM.Cells.Item(107, 13).Value = -1141136.2  # M107: 1641.2 -> -1141136.2
$ws_GSM.Cells.Item(107, 14).Value = -4814.3  # N107: -4721.8182 -> -4814.3

# GSM row 132 (G=44008)
$ws_GSM.Cells.Item(132, 8).Value = 1943.1052  # H132: 1982.8108 -> 1943.1052
$ws_GSM.Cells.Item(132, 9).Value = 1500.2812  # I132: 1533.3871 -> 1500.2812
$ws_GSM.Cells.Item(132, 11).Value = 4500.8436  # K132: 4600.1613 -> 4500.8436
$ws_GSM.Cells.Item(132, 13).Value = -1970.8436  # M132: -2070.1613 -> -1970.8436

# LTW row 17 (G=3017)
$ws_LTW.Cells.Item(17, 8).Value = 0  # H17: 9999 -> 0
$ws_LTW.Cells.Item(17, 10).Value = 0  # J17: 9999 -> 0
$ws_LTW.Cells.Item(17, 12).Value = 0  # L17: 9999 -> 0
$ws_LTW.Cells.Item(17, 14).ClearContents()  # N17: -10339 -> (removed)

# LTW row 46 (G=5282)
$ws_LTW.Cells.Item(46, 8).Value = 2538.9697  # H46: 2473.8857 -> 2538.9697
$ws_LTW.Cells.Item(46, 9).Value = 2043.4783  # I46: 1992 -> 2043.4783
$ws_LTW.Cells.Item(46, 11).Value = 2043.4783  # K46: 1992 -> 2043.4783
$ws_LTW.Cells.Item(46, 13).Value = -1855.4783  # M46: -1804 -> -1855.4783

# WVR row 37 (G=3351)
$ws_WVR.Cells.Item(37, 8).Value = 0  # H37: 24999 -> 0
$ws_WVR.Cells.Item(37, 10).Value = 0  # J37: 24999 -> 0
$ws_WVR.Cells.Item(37, 12).Value = 0  # L37: 24999 -> 0
$ws_WVR.Cells.Item(37, 14).ClearContents()  # N37: -25405 -> (removed)

# WVR row 113 (G=27752)
$ws_WVR.Cells.Item(113, 8).Value = 737.1458  # H113: 749.93616 -> 737.1458
$ws_WVR.Cells.Item(113, 9).Value = 736.0833  # I113: 742.5417 -> 736.0833
$ws_WVR.Cells.Item(113, 10).Value = 738.2083  # J113: 757.65216 -> 738.2083
$ws_WVR.Cells.Item(113, 11).Value = 2208.2499  # K113: 2227.6251 -> 2208.2499
$ws_WVR.Cells.Item(113, 12).Value = 2214.6249  # L113: 2272.95648 -> 2214.6249
$ws_WVR.Cells.Item(113, 13).Value = -38.2498999999998  # M113: -57.6251000000002 -> -38.2498999999998
$ws_WVR.Cells.Item(113, 14).Value = -6554.6249  # N113: -6612.95648 -> -6554.6249

# WVR row 135 (G=42043)
$ws_WVR.Cells.Item(135, 8).Value = 52287.5  # H135: 52520 -> 52287.5
$ws_WVR.Cells.Item(135, 10).Value = 52287.5  # J135: 52520 -> 52287.5
$ws_WVR.Cells.Item(135, 12).Value = 52287.5  # L135: 52520 -> 52287.5
$ws_WVR.Cells.Item(135, 14).Value = -62427.5  # N135: -62660 -> -62427.5

# WVR row 136 (G=44031)
$ws_WVR.Cells.Item(136, 8).Value = 30610726  # H136: 33671716 -> 30610726
$ws_WVR.Cells.Item(136, 9).Value = 90910600  # I136: 111112790 -> 90910600
$ws_WVR.Cells.Item(136, 10).Value = 460787.22  # J136: 482681.62 -> 460787.22
$ws_WVR.Cells.Item(136, 11).Value = 272731800  # K136: 333338370 -> 272731800
$ws_WVR.Cells.Item(136, 12).Value = 1382361.66  # L136: 1448044.86 -> 1382361.66
$ws_WVR.Cells.Item(136, 13).Value = -272729250  # M136: -333335820 -> -272729250
$ws_WVR.Cells.Item(136, 14).Value = -1387461.66  # N136: -1453144.86 -> -1387461.66
